# -----------------------------------------------------------------------
# Applies the Salary-Slip-word.docx edits described by the diff:
#   1. "Month     :" -> "Month    :" (one fewer space before the colon)
#   2. "Dec'23"      -> "Nov'23" (month changed, apostrophe/year re-split
#                         into separate runs: "Nov" | "'" | "2" | "3")
#   3. "20000.00"    -> "30000.00" (Basic pay changed, split into
#                         "3" | "0000" | ".00")
#   4. The "H R A" row's label and amount ("10000.00") are cleared out,
#      leaving two blank paragraphs (and the amount cell's center
#      alignment is removed along with it).
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

function Split-RangeFormatting($rng) {
    # Word COM-interop here coalesces adjacent runs that share identical
    # rPr once any edit touches/abuts their boundary. Nudging a range's
    # font size away from, then back to, its original value forces the
    # engine to materialize that range as its own run again without
    # altering the rendered formatting.
    $orig = $rng.Font.Size
    $bump = $orig + 1
    $rng.Font.Size = $bump
    $rng.Font.Size = $orig
}

# --- 1 & 2: the "Month ... : Dec'23" line -------------------------------

$monthPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Month*Dec*23*") {
        $monthPara = $cand
        break
    }
}

$pStart = $monthPara.Range.Start
$pText = $monthPara.Range.Text

$decOffset = $pText.IndexOf("Dec")
$decStart = $pStart + $decOffset
$decEnd = $decStart + 3

# "Dec" -> "Nov" (merges with the trailing "'23" into one run)
$d.Range($decStart, $decEnd).Text = "Nov"

# Re-split "Nov" | "'" | "2" | "3"
$aposStart = $decEnd
$aposEnd = $aposStart + 1
Split-RangeFormatting $d.Range($aposStart, $aposEnd)

$twoStart = $aposEnd
$threeStart = $twoStart + 1
$threeEnd = $threeStart + 1
Split-RangeFormatting $d.Range($threeStart, $threeEnd)

# 5 spaces -> 4 spaces right after "Month" (the run that directly
# precedes the "gramStart" proofing mark). "Month" is always 5 chars,
# and this whitespace run is always exactly 5 chars wide in the source.
$spacesStart = $pStart + 5
$spacesEnd = $spacesStart + 5
$d.Range($spacesStart, $spacesEnd).Text = "    "
# Re-split "Month" from the (now 4-space) run.
Split-RangeFormatting $d.Range($spacesStart, $spacesStart + 4)

# --- 3: Basic pay 20000.00 -> 30000.00 ----------------------------------

$table = $d.Tables.Item(1)
$basicRow = $null
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $row = $table.Rows.Item($r)
    if ($row.Cells.Item(1).Range.Text -like "Basic*") {
        $basicRow = $row
        break
    }
}
$amountCell = $basicRow.Cells.Item(2)
$amountStart = $amountCell.Range.Start

$d.Range($amountStart, $amountStart + 1).Text = "3"
Split-RangeFormatting $d.Range($amountStart + 1, $amountStart + 5)
Split-RangeFormatting $d.Range($amountStart + 5, $amountStart + 8)

# --- 4: clear out the "H R A" row ---------------------------------------

$hraRow = $null
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $row = $table.Rows.Item($r)
    if ($row.Cells.Item(1).Range.Text -like "H R A*") {
        $hraRow = $row
        break
    }
}

$labelPara = $hraRow.Cells.Item(1).Range.Paragraphs.Item(1)
$lStart = $labelPara.Range.Start
$lEnd = $labelPara.Range.End
$d.Range($lStart, $lEnd - 1).Text = ""

$amtPara = $hraRow.Cells.Item(2).Range.Paragraphs.Item(1)
$amtPara.Range.ParagraphFormat.Alignment = 0
$aStart = $amtPara.Range.Start
$aEnd = $amtPara.Range.End
$d.Range($aStart, $aEnd - 1).Text = ""

Write-Output $d.Content.Text
